$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so that
# values such as "1.00" or "538.28" are not silently coerced into numbers
# (which would lose trailing zeros / exact formatting, or round floats).
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.154.94"
$ws.Range("E2").Value = "  -3.48%  "
$ws.Range("D3").Value = "2.510.49"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "538.28"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "144.07"
$ws.Range("E6").Value = "  -4.65%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "2.549.91"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("D12").Value = "5.55"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "2.953.93"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "23.72"
$ws.Range("E15").Value = "  -5.40%  "
$ws.Range("D16").Value = "59.048.28"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "2.533.65"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").Value = "323.11"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "5.78"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").Value = "62.07"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  -8.78%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "2.622.62"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "0.991"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").Value = "0.0₃0778"
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "158.73"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "18.61"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "4.40"
$ws.Range("E38").Value = "  -7.91%  "
$ws.Range("D39").Value = "1.62"
$ws.Range("E39").Value = "  -7.51%  "
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "305.10"
$ws.Range("E41").Value = "  -4.99%  "
$ws.Range("D42").Value = "36.88"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "0.824"
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").Value = "0.993"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").Value = "10.76"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "125.82"
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").Value = "18.80"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -3.33%  "

# Restore the original (default) cell style now that the text values are
# safely stored, so the workbook's styles.xml / cell "s" attributes are
# unaffected by the temporary text number-format.
$textRange.Style = "Normal"
